$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the species-record data between row 2 and row 3 (columns A, B, D, E, F, G, H),
# and move the "Biotop" value (AH) from row 3 to row 2.

$ws.Range("A2").Value = 111661750
$ws.Range("B2").Value = 99581
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 221317
$ws.Range("F2").Value = "Gullklöver"
$ws.Range("G2").Value = "Trifolium aureum"
$ws.Range("H2").Value = "Pollich"
$ws.Range("AH2").Value = "Vägkant"

$ws.Range("A3").Value = 111661765
$ws.Range("B3").Value = 89953
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 3884
$ws.Range("F3").Value = "Hasselticka"
$ws.Range("G3").Value = "Dichomitus campestris"
$ws.Range("H3").Value = "(Quél.) Domański & Orlicz"
$ws.Range("AH3").Value = ""
